# Daily "remaining days" rollover update.
#
# Columns: A=行号 B=店铺名称 C=地址 D=总天(total days) E=剩余(days remaining)
#          F=开始时间(cycle start date, yyyymmdd)  G/H/I=notes
#
# Business rule (one day has passed - "today" advanced to 2025-12-15):
#   newE = D - (today - F)          # one fewer day remaining
#   if newE <= 0 then the cycle expired today, so it restarts:
#       newE = D
#       newF = today (yyyymmdd)
#   otherwise F is left unchanged.
# Rows whose F value isn't a well-formed 8-digit yyyymmdd date are left
# completely untouched (can't compute a day offset for them).

function Get-JulianDay($y, $m, $d) {
    $a = [math]::Floor((14 - $m) / 12)
    $y2 = $y + 4800 - $a
    $m2 = $m + 12 * $a - 3
    $jdn = $d + [math]::Floor((153 * $m2 + 2) / 5) + 365 * $y2 + [math]::Floor($y2 / 4) - [math]::Floor($y2 / 100) + [math]::Floor($y2 / 400) - 32045
    return $jdn
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todayYear = 2025
$todayMonth = 12
$todayDay = 15
$todayJD = Get-JulianDay $todayYear $todayMonth $todayDay
$todayStr = 20251215

$lastRow = $ws.Cells.Item(1, 4).End(-4121).Row
if ($lastRow -lt 2) {
    $lastRow = 99
}

for ($row = 2; $row -le $lastRow; $row++) {
    $D = $ws.Cells.Item($row, 4).Value2
    $E = $ws.Cells.Item($row, 5).Value2
    $F = $ws.Cells.Item($row, 6).Value2

    if ($null -eq $D -or $null -eq $F) {
        continue
    }

    $Fstr = [string]([int64]$F)
    if ($Fstr.Length -ne 8) {
        # Malformed date (e.g. "202510929") - can't compute, skip this row.
        continue
    }

    $fy = [int]$Fstr.Substring(0, 4)
    $fm = [int]$Fstr.Substring(4, 2)
    $fd = [int]$Fstr.Substring(6, 2)
    $fJD = Get-JulianDay $fy $fm $fd

    $newE = $D - ($todayJD - $fJD)
    $newF = $F

    if ($newE -le 0) {
        # Cycle expired today - restart the countdown from today.
        $newE = $D
        $newF = $todayStr
    }

    if ($newE -ne $E) {
        $ws.Cells.Item($row, 5).Value = $newE
    }
    if ($newF -ne $F) {
        $ws.Cells.Item($row, 6).Value = $newF
    }
}
